$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 69
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 0.944
$ws.Range("E2").Value = 97
$ws.Range("F2").Value = 0.894
$ws.Range("G2").Value = 0.001
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 0.324
$ws.Range("J2").Value = 0.532
$ws.Range("K2").Value = 0.381
$ws.Range("L2").Value = 0.324
$ws.Range("M2").Value = 0.611
$ws.Range("N2").Value = 2.062
$ws.Range("O2").Value = 1.481
$ws.Range("P2").Value = 0.314
$ws.Range("Q2").Value = 0.541

# Row 3
$ws.Range("B3").Value = 59
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 0.823
$ws.Range("E3").Value = 99
$ws.Range("F3").Value = 0.902
$ws.Range("G3").Value = 0.001
$ws.Range("H3").Value = 5
$ws.Range("I3").Value = 0.021
$ws.Range("J3").Value = 0.161
$ws.Range("K3").Value = 0.185
$ws.Range("L3").Value = 0.021
$ws.Range("M3").Value = 0.439
$ws.Range("N3").Value = 1.934
$ws.Range("O3").Value = 0.461
$ws.Range("P3").Value = 0.223
$ws.Range("Q3").Value = 0.471

# Row 4
$ws.Range("B4").Value = 75
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 0.85
$ws.Range("E4").Value = 101
$ws.Range("F4").Value = 0.931
$ws.Range("G4").Value = 0.001
$ws.Range("H4").Value = 4
$ws.Range("I4").Value = 0.018
$ws.Range("J4").Value = 0.09
$ws.Range("K4").Value = 0.202
$ws.Range("L4").Value = 0.018
$ws.Range("M4").Value = 0.242
$ws.Range("N4").Value = 2.518
$ws.Range("O4").Value = 0.196
$ws.Range("P4").Value = 0.155
$ws.Range("Q4").Value = 0.415

# Row 5
$ws.Range("B5").Value = 141
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 0.824
$ws.Range("E5").Value = 96
$ws.Range("F5").Value = 0.898
$ws.Range("G5").Value = 0.001
$ws.Range("H5").Value = 5
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0.358
$ws.Range("N5").Value = 1.879
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0.132
$ws.Range("Q5").Value = 0.531

# Row 6
$ws.Range("B6").Value = 33
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 0.858
$ws.Range("E6").Value = 100
$ws.Range("F6").Value = 0.958
$ws.Range("G6").Value = 0.001
$ws.Range("H6").Value = 4
$ws.Range("I6").Value = 0.233
$ws.Range("J6").Value = 0.524
$ws.Range("K6").Value = 0.438
$ws.Range("L6").Value = 0.233
$ws.Range("M6").Value = 0.5
$ws.Range("N6").Value = 2.285
$ws.Range("O6").Value = 0.8110000000000001
$ws.Range("P6").Value = 0.1
$ws.Range("Q6").Value = 0.528

$text_R2 = @"
Subgraph 1 (Nodes): [62, 54, 37, 41, 40, 26, 50] - Density: 1
Subgraph 2 (Nodes): [35, 66, 18, 25, 33, 11, 36] - Density: 0.857143
Subgraph 3 (Nodes): [67, 15, 12, 59, 10, 29] - Density: 1
Subgraph 4 (Nodes): [31, 44, 6, 42, 16, 22, 45, 32, 7] - Density: 0.916647
----------------------------------------------------
"@
$ws.Range("R2").Value = $text_R2

$text_R3 = @"
Subgraph 1 (Nodes): [25, 4, 3, 33, 5, 45] - Density: 1
Subgraph 2 (Nodes): [46, 9, 20, 43, 47, 44, 56, 1, 37] - Density: 0.845224
Subgraph 3 (Nodes): [21, 51, 0, 42, 26, 16, 2, 50, 19, 8, 39, 11, 32, 22] - Density: 0.813187
----------------------------------------------------
"@
$ws.Range("R3").Value = $text_R3

$text_R4 = @"
Subgraph 1 (Nodes): [39, 70, 46, 35, 31, 23, 64, 24, 45, 26, 30, 34, 20, 14] - Density: 0.837919
Subgraph 2 (Nodes): [28, 18, 42, 53, 15, 40, 32, 43, 2, 11, 50, 74, 3, 59] - Density: 0.843412
Subgraph 3 (Nodes): [10, 17, 44, 60, 5, 38, 41, 63, 25, 72, 36, 8, 67] - Density: 0.853144
Subgraph 4 (Nodes): [68, 47, 6, 0, 19, 52, 4, 21, 51, 12] - Density: 0.816659
Subgraph 5 (Nodes): [1, 37, 58, 54, 55, 49, 69, 56, 22, 13] - Density: 0.808335
----------------------------------------------------
"@
$ws.Range("R4").Value = $text_R4

$text_R5 = @"
Subgraph 1 (Nodes): [3, 64, 84, 119, 12, 99, 94, 69, 55, 24, 46, 76, 98, 111, 60, 86, 53, 88, 48, 79, 0, 66, 116, 39, 92, 74, 20, 108, 67, 72, 27, 114, 133, 77, 51, 9, 49, 85, 56, 41, 128, 137] - Density: 0.806918
Subgraph 2 (Nodes): [17, 130, 91, 58, 96, 131, 106, 80, 59, 134, 135, 62, 21, 26] - Density: 0.813187
Subgraph 3 (Nodes): [125, 87, 127, 93, 107, 117, 44, 118, 36, 1, 83, 113, 23, 140, 105, 30, 109, 5, 95, 112, 18, 8, 11, 19, 29, 70, 138, 2, 16, 78] - Density: 0.816823
----------------------------------------------------
"@
$ws.Range("R5").Value = $text_R5

$text_R6 = @"
Subgraph 1 (Nodes): [1, 6, 13, 2, 8] - Density: 1
Subgraph 2 (Nodes): [28, 12, 21, 16, 7] - Density: 1
Subgraph 3 (Nodes): [25, 18, 22] - Density: 1
Subgraph 4 (Nodes): [29, 19, 27] - Density: 1
Subgraph 5 (Nodes): [14, 0, 30, 20] - Density: 1
----------------------------------------------------
"@
$ws.Range("R6").Value = $text_R6

$text_S2 = @"
--- DEBUG: Seeds Loaded ---
Total seeds = 2
Seed 1: { 26 37 40 41 50 54 62 } | Triangles: 35 | Density: 1
Seed 2: { 10 12 15 29 59 67 } | Triangles: 20 | Density: 1
Subgraph 1:z { 1 26 37 40 41 48 50 54 62 } N: 9 Triangles: 35 Density: 0.416647
Subgraph 2:z { 0 6 7 10 12 14 15 16 22 27 29 30 31 32 34 42 44 45 59 63 67 } N: 21 Triangles: 97 Density: 0.0729225
"@
$ws.Range("S2").Value = $text_S2

$text_S3 = @"
--- DEBUG: Seeds Loaded ---
Total seeds = 2
Seed 1: { 1 9 20 37 43 44 46 47 56 } | Triangles: 71 | Density: 0.845224
Seed 2: { 3 4 5 25 33 45 } | Triangles: 20 | Density: 1
Subgraph 1:z { 0 1 2 3 4 5 8 9 11 16 19 20 21 22 23 24 25 26 27 30 32 33 36 37 39 42 43 44 45 46 47 48 50 51 54 56 57 } N: 37 Triangles: 387 Density: 0.0497287
Subgraph 2:z { 0 1 2 3 4 5 8 9 11 16 19 20 21 22 23 24 25 26 27 30 32 33 36 37 39 42 43 44 45 46 47 48 50 51 54 56 57 } N: 37 Triangles: 387 Density: 0.0497287
"@
$ws.Range("S3").Value = $text_S3

$text_S4 = @"
--- DEBUG: Seeds Loaded ---
Total seeds = 1
Seed 1: { 5 8 10 17 25 36 38 41 44 60 63 67 72 } | Triangles: 244 | Density: 0.853144
Subgraph 1:z { 0 1 2 3 4 5 6 8 10 11 12 13 14 15 17 18 19 20 21 22 23 24 25 26 28 29 30 31 32 33 34 35 36 37 38 39 40 41 42 43 44 45 46 47 49 50 51 52 53 54 55 56 58 59 60 63 64 66 67 68 69 70 72 74 } N: 64 Triangles: 1051 Density: 0.0248179
"@
$ws.Range("S4").Value = $text_S4

$text_S5 = @"
--- DEBUG: Seeds Loaded ---
Total seeds = 0
"@
$ws.Range("S5").Value = $text_S5

$text_S6 = @"
--- DEBUG: Seeds Loaded ---
Total seeds = 3
Seed 1: { 7 12 16 21 28 } | Triangles: 10 | Density: 1
Seed 2: { 0 14 20 30 } | Triangles: 4 | Density: 1
Seed 3: { 1 2 6 8 13 } | Triangles: 10 | Density: 1
Subgraph 1:z { 0 7 12 14 16 18 19 20 21 22 25 27 28 29 30 } N: 15 Triangles: 16 Density: 0.035166
Subgraph 2:z { 0 7 12 14 16 18 19 20 21 22 25 27 28 29 30 } N: 15 Triangles: 16 Density: 0.035166
Subgraph 3:z { 1 2 6 8 13 } N: 5 Triangles: 10 Density: 1
"@
$ws.Range("S6").Value = $text_S6
